$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("x_vs_t")

# Update the B-column values (X0 growth curve values) for rows 2-17
$ws.Cells.Item(2, 2).Value = 10
$ws.Cells.Item(3, 2).Value = 12.32
$ws.Cells.Item(4, 2).Value = 17.14
$ws.Cells.Item(5, 2).Value = 21.07
$ws.Cells.Item(6, 2).Value = 21.87
$ws.Cells.Item(7, 2).Value = 22.67
$ws.Cells.Item(8, 2).Value = 23.92
$ws.Cells.Item(9, 2).Value = 22.71
$ws.Cells.Item(10, 2).Value = 22.25
$ws.Cells.Item(11, 2).Value = 22.32
$ws.Cells.Item(12, 2).Value = 21.4
$ws.Cells.Item(13, 2).Value = 20.03
$ws.Cells.Item(14, 2).Value = 20.35
$ws.Cells.Item(15, 2).Value = 18.75
$ws.Cells.Item(16, 2).Value = 16.07
$ws.Cells.Item(17, 2).Value = 15.53

# Remove the now-unused trailing rows (18-22), shrinking the sheet to A1:B17
$ws.Range("A18:B22").ClearContents() | Out-Null

# Make x_vs_t the active sheet/tab with the D7 cell selected
$ws.Activate() | Out-Null
$ws.Range("D7").Select() | Out-Null
